# Update for 14 April
# Adds a new "4/13/20" deaths column (AF) to the US states deaths sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing data column (AE) into the
# new column (AF) for the header row and all 53 data rows (rows 1-54).
$ws.Range("AE1:AE54").Copy()
$ws.Range("AF1:AF54").PasteSpecial(-4122)

# New header label for column AF.
$ws.Range("AF1").Value = " 4/13/20"

# New deaths-per-state values for 4/13/20, in row order (row 2 = Alabama ... row 54 = Wyoming).
$values = @(103, 8, 122, 30, 731, 304, 602, 41, 52, 499, 480, 5, 9, 33, 794, 350, 43, 62, 104, 884, 19, 262, 844, 1602, 70, 98, 122, 7, 18, 120, 23, 2443, 31, 10056, 107, 8, 274, 99, 53, 590, 45, 73, 87, 6, 109, 295, 18, 28, 149, 523, 9, 154, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 32).Value = $values[$i]
}

# Update the active selection to reflect the newly added column, matching
# the author's last selection before saving.
$ws.Range("AF2").Select()
